$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The scraper added two new stat columns ("height" and "weight") right after
# "fumbles". Inserting at E:F shifts the existing "fantasy points" column
# (and its formatting/values) from E to G automatically.
$ws.Range("E:F").Insert()

# New column headers for the inserted columns.
$ws.Range("E1").Value = "height"
$ws.Range("F1").Value = "weight"

# New column data: height is a constant 6.5 and weight is a constant 255 for
# every player-week row.
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 5).Value = 6.5
    $ws.Cells.Item($r, 6).Value = 255
}

$ws.Range("A1").Select()
